$wb = $excel.ActiveWorkbook

# --- Sheet2 (SitewideSearchEs) ---
$ws2 = $wb.Worksheets.Item("SitewideSearchEs")

# Row 3: drop "tabaco" (B); pull row 4's B/C up into row 3.
$ws2.Range("B3").Value = "linfoma"
$ws2.Range("C3").Value = "cáncer de hígado"

# Row 4: only "macrófago" (old A4) survives, in column A.
$ws2.Range("A4").Value = "macrófago"
$ws2.Range("B4:C4").ClearContents() | Out-Null

# Row 5: only "dermis" (old A5) survives.
$ws2.Range("A5").Value = "dermis"

# Drop the now-unused trailing rows 6-8 entirely.
$ws2.Rows("6:8").Delete() | Out-Null

$ws2.Range("A6").Select() | Out-Null

# --- Sheet1 (SitewideSearch) ---
$ws1 = $wb.Worksheets.Item("SitewideSearch")

# Row 2: drop "cancer" (C) and "shark" (D, which becomes the numeric 250
# that used to live in D3) -- pull row 3's C/D up into row 2.
$ws1.Range("C2").Value = "glioma"
$ws1.Range("D2").Clear() | Out-Null
$ws1.Range("D2").Value = 250
$ws1.Range("D2").NumberFormat = "@"

# Row 3: drop "Hematologic/Blood Cancers" (B) and the Lorem-ipsum text (E);
# pull row 4's B..E up into row 3, and clear D3's leftover text style.
$ws1.Range("B3").Value = "LiveHelp"
$ws1.Range("C3").Value = "abdominoperineal resection "
$ws1.Range("D3").Clear() | Out-Null
$ws1.Range("D3").Value = "Dr. Norman E. Sharpless"
$ws1.Range("E3").Value = "10001110101 10001110101"

# Row 4: only "[F-18]HX4" (old A7) survives, in column A.
$ws1.Range("A4").Value = "[F-18]HX4"
$ws1.Range("B4:E4").ClearContents() | Out-Null

# Row 5: only "ipilimumab" (old A8) survives, in column A; E5's styled
# empty cell is untouched.
$ws1.Range("A5").Value = "ipilimumab"

# Drop the now-unused trailing rows 6-8 entirely.
$ws1.Rows("6:8").Delete() | Out-Null

$ws1.Range("A6").Select() | Out-Null
